$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format for numeric-looking price values so they are not
# auto-converted to numbers by Excel (matching source data as text).
$textCells = @("D5", "D11", "D15", "D16", "D18", "D19", "D25", "D27", "D32", "D38", "D40", "D44", "D46", "D47", "D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply updated cell values per the cryptos list refresh.
$ws.Range("D2").Value = "26.997.09"
$ws.Range("E2").Value = "  -0.22%  "
$ws.Range("E3").Value = "  +0.13%  "
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("D5").Value = "215.12"
$ws.Range("E5").Value = "  -0.50%  "
$ws.Range("E6").Value = "  +1.48%  "
$ws.Range("E7").Value = "  +0.08%  "
$ws.Range("E8").Value = "  +0.20%  "
$ws.Range("E9").Value = "  +0.21%  "
$ws.Range("E10").Value = "  +0.69%  "
$ws.Range("D11").Value = "0.0888"
$ws.Range("E11").Value = "  -0.25%  "
$ws.Range("D12").Value = "1.912.60"
$ws.Range("E12").Value = "  +0.13%  "
$ws.Range("D13").Value = "1.692.16"
$ws.Range("E13").Value = "  +1.01%  "
$ws.Range("D15").Value = "0.530"
$ws.Range("E15").Value = "  +1.59%  "
$ws.Range("D16").Value = "65.79"
$ws.Range("E16").Value = "  -0.03%  "
$ws.Range("D17").Value = "27.002.66"
$ws.Range("E17").Value = "  -0.25%  "
$ws.Range("D18").Value = "8.17"
$ws.Range("E18").Value = "  +5.89%  "
$ws.Range("D19").Value = "237.21"
$ws.Range("E19").Value = "  +0.79%  "
$ws.Range("D20").Value = "0.0₃0735"
$ws.Range("E20").Value = "  -0.49%  "
$ws.Range("E21").Value = "  +0.08%  "
$ws.Range("E22").Value = "  -0.54%  "
$ws.Range("E23").Value = "  -0.76%  "
$ws.Range("E24").Value = "  -2.12%  "
$ws.Range("D25").Value = "146.24"
$ws.Range("E25").Value = "  +0.63%  "
$ws.Range("E26").Value = "  +0.77%  "
$ws.Range("D27").Value = "16.13"
$ws.Range("E27").Value = "  +1.28%  "
$ws.Range("E28").Value = "  -1.49%  "
$ws.Range("E29").Value = "  +0.23%  "
$ws.Range("E30").Value = "  +0.01%  "
$ws.Range("E31").Value = "  -0.39%  "
$ws.Range("D32").Value = "3.32"
$ws.Range("E32").Value = "  -0.14%  "
$ws.Range("D33").Value = "1.478.30"
$ws.Range("E33").Value = "  +1.31%  "
$ws.Range("E35").Value = "  +4.96%  "
$ws.Range("E37").Value = "  +2.43%  "
$ws.Range("D38").Value = "0.0174"
$ws.Range("E38").Value = "  +2.42%  "
$ws.Range("E39").Value = "  +1.18%  "
$ws.Range("D40").Value = "5.86"
$ws.Range("E40").Value = "  -3.45%  "
$ws.Range("E41").Value = "  +1.11%  "
$ws.Range("E43").Value = "  +1.93%  "
$ws.Range("D44").Value = "67.47"
$ws.Range("E44").Value = "  +2.29%  "
$ws.Range("D45").Value = "1.818.95"
$ws.Range("E45").Value = "  -0.05%  "
$ws.Range("D46").Value = "0.782"
$ws.Range("E46").Value = "  +0.07%  "
$ws.Range("D47").Value = "90.46"
$ws.Range("E47").Value = "  +0.19%  "
$ws.Range("E48").Value = "  +0.89%  "
$ws.Range("E49").Value = "  -0.43%  "
$ws.Range("E50").Value = "  +1.36%  "
$ws.Range("B51").Value = "Cronos"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D51").Value = "0.0508"
$ws.Range("E51").Value = "  -0.23%  "
